# Update the multiplication-fact answers in the table.
# Each Find/Replace targets the full "A×B=C" text of one run; order matters
# because the new value "202×3=606" (introduced later below) must not be
# re-caught by the earlier rule that replaces the original "202×3=606".
$d = $word.ActiveDocument

$d.Content.Find.Execute("514×6=3084", $true, $false, $false, $false, $false, $true, 1, $false, "985×4=3940", 2) | Out-Null
$d.Content.Find.Execute("846×5=4230", $true, $false, $false, $false, $false, $true, 1, $false, "299×4=1196", 2) | Out-Null
$d.Content.Find.Execute("169×8=1352", $true, $false, $false, $false, $false, $true, 1, $false, "369×7=2583", 2) | Out-Null
$d.Content.Find.Execute("123×6=738", $true, $false, $false, $false, $false, $true, 1, $false, "937×4=3748", 2) | Out-Null
$d.Content.Find.Execute("202×3=606", $true, $false, $false, $false, $false, $true, 1, $false, "831×7=5817", 2) | Out-Null
$d.Content.Find.Execute("544×3=1632", $true, $false, $false, $false, $false, $true, 1, $false, "515×5=2575", 2) | Out-Null
$d.Content.Find.Execute("975×7=6825", $true, $false, $false, $false, $false, $true, 1, $false, "402×9=3618", 2) | Out-Null
$d.Content.Find.Execute("411×2=822", $true, $false, $false, $false, $false, $true, 1, $false, "321×7=2247", 2) | Out-Null
$d.Content.Find.Execute("474×8=3792", $true, $false, $false, $false, $false, $true, 1, $false, "965×3=2895", 2) | Out-Null
$d.Content.Find.Execute("930×3=2790", $true, $false, $false, $false, $false, $true, 1, $false, "868×8=6944", 2) | Out-Null
$d.Content.Find.Execute("939×4=3756", $true, $false, $false, $false, $false, $true, 1, $false, "916×7=6412", 2) | Out-Null
$d.Content.Find.Execute("114×8=912", $true, $false, $false, $false, $false, $true, 1, $false, "718×7=5026", 2) | Out-Null
$d.Content.Find.Execute("144×8=1152", $true, $false, $false, $false, $false, $true, 1, $false, "110×3=330", 2) | Out-Null
$d.Content.Find.Execute("570×6=3420", $true, $false, $false, $false, $false, $true, 1, $false, "437×9=3933", 2) | Out-Null
$d.Content.Find.Execute("325×4=1300", $true, $false, $false, $false, $false, $true, 1, $false, "818×9=7362", 2) | Out-Null
$d.Content.Find.Execute("788×8=6304", $true, $false, $false, $false, $false, $true, 1, $false, "721×3=2163", 2) | Out-Null
$d.Content.Find.Execute("895×6=5370", $true, $false, $false, $false, $false, $true, 1, $false, "310×9=2790", 2) | Out-Null
$d.Content.Find.Execute("130×3=390", $true, $false, $false, $false, $false, $true, 1, $false, "606×6=3636", 2) | Out-Null
$d.Content.Find.Execute("647×5=3235", $true, $false, $false, $false, $false, $true, 1, $false, "466×6=2796", 2) | Out-Null
$d.Content.Find.Execute("340×9=3060", $true, $false, $false, $false, $false, $true, 1, $false, "202×3=606", 2) | Out-Null
$d.Content.Find.Execute("944×8=7552", $true, $false, $false, $false, $false, $true, 1, $false, "575×6=3450", 2) | Out-Null
$d.Content.Find.Execute("250×6=1500", $true, $false, $false, $false, $false, $true, 1, $false, "555×6=3330", 2) | Out-Null
$d.Content.Find.Execute("757×2=1514", $true, $false, $false, $false, $false, $true, 1, $false, "668×9=6012", 2) | Out-Null
$d.Content.Find.Execute("556×7=3892", $true, $false, $false, $false, $false, $true, 1, $false, "679×2=1358", 2) | Out-Null
$d.Content.Find.Execute("167×9=1503", $true, $false, $false, $false, $false, $true, 1, $false, "576×6=3456", 2) | Out-Null
